# edit.ps1
# Updates the "想去人数" (interest-count) column F on each of the four
# worksheets (展览 / 演出 / 本地生活 / 全部类型) of the 上海-漫展信息
# workbook to reflect newly generated numbers, per commit
# "Update gh-pages to output generated at 456a3b4".
#
# Only column F (numeric "want-to-go" counters) values change; everything
# else in the workbook stays as-is.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 120   # F2: 117 -> 120
$ws.Cells.Item(3, 6).Value = 947   # F3: 945 -> 947
$ws.Cells.Item(4, 6).Value = 590   # F4: 588 -> 590
$ws.Cells.Item(5, 6).Value = 2797   # F5: 2772 -> 2797
$ws.Cells.Item(6, 6).Value = 769   # F6: 766 -> 769
$ws.Cells.Item(7, 6).Value = 586   # F7: 585 -> 586
$ws.Cells.Item(8, 6).Value = 586   # F8: 585 -> 586
$ws.Cells.Item(9, 6).Value = 71   # F9: 67 -> 71
$ws.Cells.Item(10, 6).Value = 653   # F10: 652 -> 653
$ws.Cells.Item(11, 6).Value = 372   # F11: 370 -> 372
$ws.Cells.Item(12, 6).Value = 415   # F12: 404 -> 415
$ws.Cells.Item(14, 6).Value = 2141   # F14: 2140 -> 2141
$ws.Cells.Item(15, 6).Value = 1244   # F15: 1242 -> 1244
$ws.Cells.Item(16, 6).Value = 727   # F16: 724 -> 727
$ws.Cells.Item(17, 6).Value = 14   # F17: 13 -> 14
$ws.Cells.Item(18, 6).Value = 2649   # F18: 2646 -> 2649
$ws.Cells.Item(19, 6).Value = 10   # F19: 9 -> 10
$ws.Cells.Item(20, 6).Value = 37   # F20: 36 -> 37
$ws.Cells.Item(24, 6).Value = 555   # F24: 541 -> 555
$ws.Cells.Item(27, 6).Value = 547   # F27: 543 -> 547
$ws.Cells.Item(28, 6).Value = 559   # F28: 557 -> 559
$ws.Cells.Item(31, 6).Value = 365   # F31: 363 -> 365
$ws.Cells.Item(32, 6).Value = 4634   # F32: 4630 -> 4634
$ws.Cells.Item(33, 6).Value = 207   # F33: 201 -> 207

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(7, 6).Value = 324   # F7: 323 -> 324
$ws.Cells.Item(8, 6).Value = 343   # F8: 341 -> 343
$ws.Cells.Item(12, 6).Value = 161   # F12: 160 -> 161
$ws.Cells.Item(19, 6).Value = 1767   # F19: 1766 -> 1767
$ws.Cells.Item(23, 6).Value = 291   # F23: 289 -> 291
$ws.Cells.Item(32, 6).Value = 494   # F32: 492 -> 494
$ws.Cells.Item(33, 6).Value = 11   # F33: 9 -> 11

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 1444   # F4: 1442 -> 1444
$ws.Cells.Item(6, 6).Value = 564   # F6: 563 -> 564
$ws.Cells.Item(7, 6).Value = 196   # F7: 194 -> 196
$ws.Cells.Item(8, 6).Value = 227   # F8: 226 -> 227

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 1444   # F3: 1442 -> 1444
$ws.Cells.Item(4, 6).Value = 564   # F4: 563 -> 564
$ws.Cells.Item(5, 6).Value = 120   # F5: 117 -> 120
$ws.Cells.Item(6, 6).Value = 196   # F6: 194 -> 196
$ws.Cells.Item(8, 6).Value = 947   # F8: 945 -> 947
$ws.Cells.Item(9, 6).Value = 590   # F9: 588 -> 590
$ws.Cells.Item(10, 6).Value = 2797   # F10: 2773 -> 2797
$ws.Cells.Item(11, 6).Value = 769   # F11: 766 -> 769
$ws.Cells.Item(12, 6).Value = 586   # F12: 585 -> 586
$ws.Cells.Item(13, 6).Value = 586   # F13: 585 -> 586
$ws.Cells.Item(14, 6).Value = 72   # F14: 67 -> 72
$ws.Cells.Item(15, 6).Value = 653   # F15: 652 -> 653
$ws.Cells.Item(17, 6).Value = 415   # F17: 404 -> 415
$ws.Cells.Item(18, 6).Value = 324   # F18: 323 -> 324
$ws.Cells.Item(19, 6).Value = 343   # F19: 341 -> 343
$ws.Cells.Item(22, 6).Value = 2141   # F22: 2140 -> 2141
$ws.Cells.Item(23, 6).Value = 1244   # F23: 1242 -> 1244
$ws.Cells.Item(24, 6).Value = 727   # F24: 724 -> 727
$ws.Cells.Item(26, 6).Value = 14   # F26: 13 -> 14
$ws.Cells.Item(27, 6).Value = 2650   # F27: 2646 -> 2650
$ws.Cells.Item(28, 6).Value = 10   # F28: 9 -> 10
$ws.Cells.Item(30, 6).Value = 37   # F30: 36 -> 37
$ws.Cells.Item(35, 6).Value = 227   # F35: 226 -> 227
$ws.Cells.Item(37, 6).Value = 555   # F37: 542 -> 555
$ws.Cells.Item(38, 6).Value = 555   # F38: 542 -> 555
$ws.Cells.Item(40, 6).Value = 547   # F40: 543 -> 547
$ws.Cells.Item(41, 6).Value = 559   # F41: 557 -> 559
$ws.Cells.Item(42, 6).Value = 291   # F42: 289 -> 291
$ws.Cells.Item(45, 6).Value = 365   # F45: 363 -> 365
$ws.Cells.Item(47, 6).Value = 4634   # F47: 4630 -> 4634
$ws.Cells.Item(48, 6).Value = 207   # F48: 201 -> 207
$ws.Cells.Item(50, 6).Value = 494   # F50: 492 -> 494

